# Update the "Metadata" sheet (active sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Version bump: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date bump
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; now "Alvearie Team"
$ws.Range("B9").Value = "Alvearie Team"

# The "Contact" row becomes a "Jurisdiction" row
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -
# remove it entirely, shifting everything below up by one row.
$ws.Rows.Item(11).Delete()

# Update the "Elements" sheet: the Extension row's Short/Definition columns
# now describe the Split Method extension instead of the generic placeholder.
$ws2 = $wb.Worksheets.Item("Elements")
$ws2.Range("K2").Value = "Split Method"
$ws2.Range("L2").Value = "Method used to identify the matched resource to split"
